# Commit: "adding the template as an internal data file"
#
# - Rename the sheets to their internal/data-file identifiers.
# - Move the active tab / selection from the "location" (formerly
#   "Locations") sheet to the "event" (formerly "Events") sheet.

$wb = $excel.ActiveWorkbook

$wsLocation = $wb.Worksheets.Item(1)   # was "Locations" -> rId1 -> sheet1.xml
$wsEvent    = $wb.Worksheets.Item(2)   # was "Events"    -> rId2 -> sheet2.xml

# Rename worksheets.
$wsLocation.Name = "location"
$wsEvent.Name    = "event"

# Update the selection remembered on the (no-longer-active) location sheet.
$wsLocation.Range("A25").Select()

# Make the event sheet the active tab, with its own remembered selection.
$wsEvent.Activate()
$wsEvent.Range("B15").Select()
